# ---------------------------------------------------------------------------
# Applies the edit described by the diff:
#  1. In "ODI Batting", the empty INNING_NUMBER cells (rows where the player
#     did not bat) are cleared out so the cell no longer exists.
#  2. A new worksheet "ODI Batting Extra" is appended after "ODI Bowling"
#     with MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
#     PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH columns.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$battingWs = $wb.Worksheets.Item("ODI Batting")

# ---------------------------------------------------------------------------
# 1. Clear the now-empty B (INNING_NUMBER) cells on "ODI Batting" for the
#    matches where the player did not bat at all.
# ---------------------------------------------------------------------------
$emptyInningCells = @("B2", "B3", "B6", "B7", "B9", "B14", "B16")
foreach ($addr in $emptyInningCells) {
    $battingWs.Range($addr).Value = ""
}

# ---------------------------------------------------------------------------
# 2. Create the new "ODI Batting Extra" worksheet as the last sheet in the
#    workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extraWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$extraWs.Name = "ODI Batting Extra"

# Preserve the outline / page setup properties used by the other sheets.
$extraWs.Outline.SummaryRow = 1
$extraWs.Outline.SummaryColumn = 1

# Copy the header formatting (bold, bordered, centered) used on the other
# sheets and then set the new header text.
$battingWs.Range("A1:F1").Copy()
$extraWs.Range("A1:F1").PasteSpecial(-4122)   # xlPasteFormats

$extraWs.Range("A1").Value = "MATCH_CODE"
$extraWs.Range("B1").Value = "BATTING_POSITION"
$extraWs.Range("C1").Value = "NUM_4"
$extraWs.Range("D1").Value = "NUM_6"
$extraWs.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extraWs.Range("F1").Value = "MAN_OF_MATCH"

# Helper: write a value as literal text (keeps numeric-looking strings such
# as "1" or "10.13%" from being reinterpreted as numbers/percentages).  An
# empty string leaves a present-but-blank cell behind instead of deleting it.
function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    if ($val -eq "") {
        $cell.Value = $null
    } else {
        $cell.Value = $val
    }
}

# Helper: write a value as a real number, leaving a present-but-blank cell
# behind (instead of deleting it) when the value is blank.
function Set-NumberCell($cell, $val) {
    if ($val -eq "") {
        $cell.NumberFormat = "@"
        $cell.Value = $null
    } else {
        $cell.Value = [double]$val
    }
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$rows = @(
    @("4098", "",  "",  "",  "",       "NO"),
    @("4533", "",  "",  "",  "",       "NO"),
    @("4535", "6", "1", "0", "10.13%", "NO"),
    @("4536", "7", "2", "1", "12.45%", "NO"),
    @("4657", "7", "",  "",  "",       "NO"),
    @("4658", "7", "",  "",  "",       "NO"),
    @("4669", "7", "3", "3", "12.09%", "NO"),
    @("4673", "",  "",  "",  "",       "NO"),
    @("4676", "7", "5", "1", "23.29%", "NO"),
    @("4679", "6", "0", "0", "10.22%", "NO"),
    @("4682", "4", "1", "0", "4.14%",  "NO"),
    @("4685", "",  "",  "",  "",       "NO"),
    @("4691", "",  "",  "",  "",       "NO"),
    @("4692", "",  "",  "",  "",       "NO"),
    @("4695", "7", "",  "",  "",       "NO"),
    @("4697", "7", "1", "0", "2.34%",  "NO")
)

$r = 2
foreach ($row in $rows) {
    Set-TextCell   $extraWs.Cells.Item($r, 1) $row[0]   # MATCH_CODE
    Set-NumberCell $extraWs.Cells.Item($r, 2) $row[1]   # BATTING_POSITION
    Set-TextCell   $extraWs.Cells.Item($r, 3) $row[2]   # NUM_4
    Set-TextCell   $extraWs.Cells.Item($r, 4) $row[3]   # NUM_6
    Set-TextCell   $extraWs.Cells.Item($r, 5) $row[4]   # PERCENT_RUNS_OF_TOTAL
    Set-TextCell   $extraWs.Cells.Item($r, 6) $row[5]   # MAN_OF_MATCH
    $r++
}

# Restore the originally active sheet so the workbook view state is
# unchanged (the new sheet is simply appended, not selected).
$wb.Worksheets.Item("Player Info").Activate()
